$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 39 ("Robot LEDs") is removed wholesale so we can rebuild the
#     bottom of the checklist: a fresh "Robot Vision" row, followed by a
#     bare "Robot LEDs" row that only carries a class name. ---
$ws.Rows.Item(39).Delete()

# New row 39: the old row-38 content ("Robot Vision" / X / formula) moves
# down here, with Notes (E:F) merged but left blank.
$ws.Range("B39").Value = "Robot Vision"
$ws.Range("C39").Value = "X"
$ws.Range("D39").Formula = '=IF(EXACT(LOWER(C39), "check"), "ü", "û")'
$ws.Range("E39:F39").Merge()

# New row 40: just the class name, nothing else.
$ws.Range("B40").Value = "Robot LEDs"

# Row 38 keeps its original C/D (X / formula), only the class name and the
# note change: the robot now tests the experimental NetworkTables client.
$ws.Range("B38").Value = "Expiremental NetTablesClient"
$ws.Range("E38").Value = "don't run customOutputStream"

# The active selection after this edit sits on the (now blank) E38:F38
# merged note cell.
$ws.Range("E38").Select()
